# Generate Report for Handback
#
# The localization-status report is refreshed once handback has completed:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     on the Overview sheet and on each per-language detail sheet.
#   - Each per-language sheet gains "Latest Target File" / "Latest Handback File"
#     hyperlink columns (F/G) pointing at the source markdown / translated xlf.
#   - "Latest Handback DateTime" (H) is stamped with the actual handback time
#     instead of the zero date.

$wb = $excel.ActiveWorkbook

$statusNew = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns for both rows ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusNew
$overview.Range("C2").Value = $statusNew
$overview.Range("B3").Value = $statusNew
$overview.Range("C3").Value = $statusNew

function Update-LanguageSheet {
    param($Workbook, $SheetName, $XlfFileName, $HandbackDateTime, $StatusNew)

    $ws = $Workbook.Worksheets.Item($SheetName)

    $mdFileName = "320188b0-d72c-41c3-9bb7-f41b58c0e308.md"
    $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/84d1674343a6848ab19ee93abb3c33cf4c4ae844/e2e/$mdFileName"
    $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad551b47c84f3535d6f3a50af7a587d1d8fae4a9/ol-handoff/OpenLocalizationTestOrg/oltest.$SheetName/ci/ht/$XlfFileName"

    # Status column (Status = "Ready for handoff" -> "Handed back: in sync with en-US")
    $ws.Range("C2").Value = $StatusNew
    $ws.Range("C3").Value = $StatusNew

    # Latest Target File (F) / Latest Handback File (G) hyperlinks, rows 2 & 3
    $ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl, "", "", $mdFileName)
    $ws.Hyperlinks.Add($ws.Range("G2"), $xlfUrl, "", "", $XlfFileName)
    $ws.Hyperlinks.Add($ws.Range("F3"), $mdUrl, "", "", $mdFileName)
    $ws.Hyperlinks.Add($ws.Range("G3"), $xlfUrl, "", "", $XlfFileName)

    $ws.Range("F2").Style = "Hyperlink"
    $ws.Range("G2").Style = "Hyperlink"
    $ws.Range("F3").Style = "Hyperlink"
    $ws.Range("G3").Style = "Hyperlink"

    # Latest Handback DateTime (H) now carries the real handback timestamp
    $ws.Range("H2").Value = $HandbackDateTime
    $ws.Range("H3").Value = $HandbackDateTime
}

Update-LanguageSheet $wb "zh-cn" "320188b0-d72c-41c3-9bb7-f41b58c0e308.d86f4dfb38f8d79df51a423d2d5f3f6c9a795df5.zh-cn.xlf" "2016-03-25 07:59:19" $statusNew
Update-LanguageSheet $wb "de-de" "320188b0-d72c-41c3-9bb7-f41b58c0e308.d86f4dfb38f8d79df51a423d2d5f3f6c9a795df5.de-de.xlf" "2016-03-25 07:59:26" $statusNew
